# Applies the "Progetto-html" re-test edit described by the commit:
# "Testato progetto HTML con locatori generati con il tool per robula+"
#
# - Updates the report title shared string.
# - Rewrites the two result tables (rows 4-9 "Analitica", rows 11-16
#   "Totali") with the new run's raw counts, turning the old D-column
#   formula (C-E-F) into a literal value while keeping the G-column
#   percentage formula.
# - Applies center+middle alignment to the rewritten C:F cells.
# - Updates the two totals rows (20/21).
# - Restores the saved selection to G22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title text -------------------------------------------------------
$ws.Range("A1").Value = "Applicazione: Progetto-html - Confronto Robustezza Locatori"

# --- Data tables --------------------------------------------------------
# Each row: C (total), D (success count, now a literal), E (fragility
# failures), F (obsolescence failures). G keeps its %-formula and just
# recalculates automatically once C/E change.
$rows = @(
    @{ R = 4;  C = 40; D = 38; E = 1;  F = 1 },
    @{ R = 5;  C = 40; D = 28; E = 11; F = 1 },
    @{ R = 6;  C = 40; D = 30; E = 9;  F = 1 },
    @{ R = 7;  C = 40; D = 38; E = 1;  F = 1 },
    @{ R = 8;  C = 40; D = 36; E = 3;  F = 1 },
    @{ R = 9;  C = 40; D = 35; E = 4;  F = 1 },
    @{ R = 11; C = 40; D = 34; E = 4;  F = 2 },
    @{ R = 12; C = 40; D = 28; E = 10; F = 2 },
    @{ R = 13; C = 40; D = 30; E = 8;  F = 2 },
    @{ R = 14; C = 40; D = 38; E = 0;  F = 2 },
    @{ R = 15; C = 40; D = 35; E = 3;  F = 2 },
    @{ R = 16; C = 40; D = 35; E = 3;  F = 2 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 3).Value = $row.C   # C
    $ws.Cells.Item($r, 4).Value = $row.D   # D (now a literal, formula dropped)
    $ws.Cells.Item($r, 5).Value = $row.E   # E
    $ws.Cells.Item($r, 6).Value = $row.F   # F

    foreach ($col in 3..6) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.HorizontalAlignment = -4108   # xlHAlignCenter / xlCenter
        $cell.VerticalAlignment = -4108     # xlVAlignCenter / xlCenter
    }
}

# --- Totals rows --------------------------------------------------------
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 1

$ws.Range("E21").Value = 0

# --- Selection ------------------------------------------------------------
[void]$ws.Range("G22").Select()
